$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        [string]$Addr,
        [string]$Val
    )
    $c = $ws.Range($Addr)
    $c.Value = "'" + $Val
    $c.Style = "Normal"
}

Set-CellText "D2" '64.155.77'
Set-CellText "E2" '  -0.20%  '
Set-CellText "D3" '3.479.51'
Set-CellText "D4" '0.999'
Set-CellText "E4" '  -0.04%  '
Set-CellText "D5" '585.37'
Set-CellText "E5" '  -0.15%  '
Set-CellText "D6" '131.97'
Set-CellText "E6" '  -1.73%  '
Set-CellText "E7" '  +0.02%  '
Set-CellText "E8" '  -1.07%  '
Set-CellText "D9" '7.71'
Set-CellText "E9" '  +6.12%  '
Set-CellText "E10" '  -1.16%  '
Set-CellText "E11" '  -0.08%  '
Set-CellText "D12" '4.071.81'
Set-CellText "E12" '  -0.41%  '
Set-CellText "E13" '  +0.03%  '
Set-CellText "E14" '  -2.51%  '
Set-CellText "D15" '3.476.98'
Set-CellText "E15" '  -0.46%  '
Set-CellText "D16" '64.114.82'
Set-CellText "E16" '  -0.36%  '
Set-CellText "D17" '25.08'
Set-CellText "E17" '  -2.80%  '
Set-CellText "D18" '9.99'
Set-CellText "E18" '  +1.02%  '
Set-CellText "E19" '  -1.23%  '
Set-CellText "E20" '  -1.72%  '
Set-CellText "D21" '385.03'
Set-CellText "E21" '  -2.66%  '
Set-CellText "D23" '3.617.05'
Set-CellText "E23" '  -0.50%  '
Set-CellText "D24" '74.55'
Set-CellText "E24" '  +0.20%  '
Set-CellText "B26" 'PEPE'
Set-CellText "C26" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-CellText "D26" '0.0000112'
Set-CellText "E26" '  -2.50%  '
Set-CellText "B27" 'Binance-PegBSC-USD'
Set-CellText "C27" 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-CellText "D27" '1.00'
Set-CellText "E27" '  -0.01%  '
Set-CellText "B28" 'PancakeSwap'
Set-CellText "C28" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-CellText "D28" '2.22'
Set-CellText "E28" '  -0.52%  '
Set-CellText "B29" 'RenderToken'
Set-CellText "C29" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-CellText "D29" '7.12'
Set-CellText "E29" '  -3.78%  '
Set-CellText "B30" 'InternetComputer(DFINITY)'
Set-CellText "C30" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-CellText "D30" '7.98'
Set-CellText "E30" '  -3.51%  '
Set-CellText "B31" 'Kaspa'
Set-CellText "C31" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-CellText "D31" '0.155'
Set-CellText "E31" '  +2.93%  '
Set-CellText "B32" 'Fetch.AI'
Set-CellText "C32" 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-CellText "D32" '1.43'
Set-CellText "E32" '  -4.28%  '
Set-CellText "B33" 'RenzoRestakedETH'
Set-CellText "C33" 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
Set-CellText "D33" '3.504.71'
Set-CellText "E33" '  -0.30%  '
Set-CellText "B34" 'USDe'
Set-CellText "C34" 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-CellText "D34" '1.00'
Set-CellText "E34" '  -0.07%  '
Set-CellText "B35" 'EthereumClassic'
Set-CellText "C35" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-CellText "D35" '22.98'
Set-CellText "E35" '  -1.91%  '
Set-CellText "B36" 'NEARProtocol'
Set-CellText "C36" 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-CellText "D36" '5.22'
Set-CellText "E36" '  +1.26%  '
Set-CellText "B37" 'Aptos'
Set-CellText "C37" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-CellText "D37" '6.77'
Set-CellText "E37" '  -1.95%  '
Set-CellText "B38" 'Monero'
Set-CellText "C38" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-CellText "D38" '162.37'
Set-CellText "E38" '  -2.23%  '
Set-CellText "E39" '  -3.46%  '
Set-CellText "B40" 'Hedera'
Set-CellText "C40" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-CellText "D40" '0.0780'
Set-CellText "E40" '  -0.41%  '
Set-CellText "B41" 'Mantle'
Set-CellText "C41" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-CellText "D41" '0.798'
Set-CellText "E41" '  -1.01%  '
Set-CellText "B42" 'FirstDigitalUSD'
Set-CellText "C42" 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-CellText "D42" '1.00'
Set-CellText "E42" '  +0.00%  '
Set-CellText "B43" 'OKB'
Set-CellText "C43" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-CellText "D43" '41.50'
Set-CellText "E43" '  -0.93%  '
Set-CellText "B44" 'Filecoin'
Set-CellText "C44" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-CellText "D44" '4.34'
Set-CellText "E44" '  -1.37%  '
Set-CellText "B45" 'Stacks'
Set-CellText "C45" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-CellText "D45" '1.62'
Set-CellText "E45" '  -2.01%  '
Set-CellText "B46" 'EnergySwap'
Set-CellText "C46" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-CellText "D46" '23.43'
Set-CellText "E46" '  -7.64%  '
Set-CellText "B47" 'ONDO'
Set-CellText "C47" 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-CellText "D47" '1.13'
Set-CellText "E47" '  -3.01%  '
Set-CellText "B48" 'Cosmos'
Set-CellText "C48" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-CellText "D48" '6.73'
Set-CellText "E48" '  -0.75%  '
Set-CellText "B49" 'SuiNetwork'
Set-CellText "C49" 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-CellText "D49" '0.901'
Set-CellText "E49" '  +0.19%  '
Set-CellText "B50" 'Maker'
Set-CellText "C50" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-CellText "D50" '2.340.96'
Set-CellText "E50" '  -4.85%  '
Set-CellText "B51" 'VeChain'
Set-CellText "C51" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-CellText "D51" '0.0254'
Set-CellText "E51" '  -2.84%  '
